# "serj homepage thumbnail set!" -- add two new bullet items under the
# "^^ GRADIENT WORK" line: a Stack Overflow reference link, followed by a
# note about positioning the "serj" image outside the flexbox div.

$d = $word.ActiveDocument

# Locate the "^^ GRADIENT WORK" paragraph robustly (search avoids the
# literal "^^" since Find treats "^" as a special-code escape).
$found = $d.Content
$null = $found.Find.Execute("GRADIENT WORK")
$gradientIndex = $found.Paragraphs.Item(1).Index
$gradientPara = $d.Paragraphs.Item($gradientIndex)

# Insert two new paragraph marks right before the GRADIENT WORK paragraph's
# own mark -- Word clones the preceding pPr (ListParagraph style + numPr
# ilvl=0/numId=1), so both new paragraphs pick up the same bullet formatting
# as the rest of this list without us having to reapply it by hand.
$breakPos = $gradientPara.Range.End - 1
$d.Range($breakPos, $breakPos).InsertParagraphAfter()

$breakPos = $gradientPara.Range.End - 1
$d.Range($breakPos, $breakPos).InsertParagraphAfter()

$linkParaIndex = $gradientIndex + 1
$textParaIndex = $gradientIndex + 2

# First new bullet: hyperlink to the Stack Overflow answer.
$linkPara = $d.Paragraphs.Item($linkParaIndex)
$url = "https://stackoverflow.com/questions/885835/position-an-image-outside-of-its-container"
$insertRange = $d.Range($linkPara.Range.Start, $linkPara.Range.Start)
$insertRange.InsertBefore($url)

$linkPara = $d.Paragraphs.Item($linkParaIndex)
$linkRange = $d.Range($linkPara.Range.Start, $linkPara.Range.Start + $url.Length)
$null = $d.Hyperlinks.Add($linkRange, $url, "", "", $url)

# Second new bullet: plain-text note.
$textPara = $d.Paragraphs.Item($textParaIndex)
$noteRange = $d.Range($textPara.Range.Start, $textPara.Range.Start)
$noteRange.InsertBefore("^^ Positioning serj outside the flexbox div LOL ")
